$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Change Type" rate-change column (I) ---
$ws.Range("I2").Value = "Change Type"
$ws.Range("I3").Value = "I"

# Give the new column the same look (number format / font / fill) as the
# neighbouring "Service" column (H), then let Excel size it to fit the
# "Change Type" header text (mirrors the bestFit width Excel computed).
$ws.Range("H2:H3").Copy() | Out-Null
$ws.Range("I2:I3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = 0
$ws.Columns.Item(9).AutoFit() | Out-Null

# Extend the title merge so the banner covers the new column too.
$ws.Range("A1:I1").Merge()

# Unify the header row's fill: every header cell now uses the same
# shaded-blue background that B2/C2/D2 already had.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null
$ws.Range("C2:I2").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = 0

# Restore the correct header text/value in I2 (PasteSpecial of formats only
# should not have touched it, but make sure it's still right).
$ws.Range("I2").Value = "Change Type"

# Put the selection where Excel left it after the edit.
$ws.Range("H5").Select() | Out-Null

Write-Host "applied"
